# Specific Station Details Section
# Adds a new company (WSC), two new locations (AB, YT) tied to the
# appropriate companies, and new asset_type/location/color rows for
# the new locations (also recoloring the existing Weir/BC row).

$wb = $excel.ActiveWorkbook

# --- Companies sheet: add WSC as a new active company -----------------
$wsCompanies = $wb.Worksheets.Item("Companies")
$wsCompanies.Range("A3").Value = "WSC"
# Copy B2 ("TRUE", stored as literal text) onto B3 so the new cell keeps
# the same text data type instead of being auto-detected as a boolean.
$wsCompanies.Range("B2").Copy($wsCompanies.Range("B3"))

# --- Locations sheet: add AB (NHS) and YT (WSC) ------------------------
$wsLocations = $wb.Worksheets.Item("Locations")
$wsLocations.Range("A3").Value = "AB"
$wsLocations.Range("B3").Value = "NHS"
$wsLocations.Range("A4").Value = "YT"
$wsLocations.Range("B4").Value = "WSC"

# --- AssetTypes sheet: recolor Weir/BC and add Cableway rows for the ---
# --- new locations -------------------------------------------------------
$wsAssetTypes = $wb.Worksheets.Item("AssetTypes")
$wsAssetTypes.Range("C3").Value = "#ff0000"
$wsAssetTypes.Range("A4").Value = "Cableway"
$wsAssetTypes.Range("B4").Value = "AB"
$wsAssetTypes.Range("C4").Value = "#3eead6"
$wsAssetTypes.Range("A5").Value = "Cableway"
$wsAssetTypes.Range("B5").Value = "YT"
$wsAssetTypes.Range("C5").Value = "#c425c8"

# Keep the originally active sheet (AssetTypes) selected, matching the
# unchanged workbook.xml bookViews/activeTab state in the source diff.
$wsAssetTypes.Activate()
